# Added Some Test Data
# Fill in the previously-empty row 3 on Sheet1 with some new test values,
# including a formula that references E7 (= O6, currently 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Added Data Here"
$ws.Range("B3").Value = "To Test"
$ws.Range("C3").Value = "If I delete it"
$ws.Range("D3").Formula = "=E7"
$ws.Range("D3").NumberFormat = "_-""$""* #,##0.00_-;\-""$""* #,##0.00_-;_-""$""* ""-""??_-;_-@_-"
